# March 24 update 3
# Adds three new columns (M: renewd, N: PlanID, O: iteration) to the
# "210_11" building sheet and fills them in for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Copy the formatting of the last existing header cell (L1, bold + border
# style) onto the three new header cells before writing their text so the
# new headers look consistent with the rest of the header row.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# --- Data rows (rows 2-21) --------------------------------------------------
$lastRow = 21
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"     # column M -> renewd
    $ws.Cells.Item($r, 14).Value = 20150274    # column N -> PlanID
    $ws.Cells.Item($r, 15).Value = 9           # column O -> iteration
}
